$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix A37 value (tiny float correction)
$ws.Range("A37").Value = 45833.45888479167

# New rows of data to append (rows 38-43)
$data = @(
    @(45834.45431180556, "EVOWHEY PROTEIN", "2Kg", "37,90€"),
    @(45853.38050893519, "EVOWHEY PROTEIN", "2Kg", "37,90€"),
    @(45853.3851083449,  "EVOWHEY PROTEIN", "2Kg", "37,90€"),
    @(45853.38682178241, "EVOWHEY PROTEIN", "2Kg", "37,90€"),
    @(45853.39294818287, "EVOWHEY PROTEIN", "2Kg", "37,90€"),
    @(45853.39424525264, "EVOWHEY PROTEIN", "2Kg", "37,90€")
)

$row = 38
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row++
}

# Match the date/time number format used in column A for the new rows
$ws.Range("A37").Copy() | Out-Null
$ws.Range("A38:A43").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

